# "Fixed up Myxicola and a few more loose ends in second review"
#
# On the Materials sheet:
#   - drop the Taxon_Local_ID column entirely
#   - drop the suborder / infraorder / superfamily columns
#   - change the scientificNameAuthorship sample value from
#     ${summary.Author} to ${summary.authority}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Remove the Taxon_Local_ID column and the suborder/infraorder/superfamily
# columns. Each one is located by its row-1 header text (exact, whole-cell
# match) right before it is deleted, so re-finding on every iteration means
# we never have to reason about how earlier deletions shifted the sheet.
$colsToDelete = @("Taxon_Local_ID", "suborder", "infraorder", "superfamily")
foreach ($colName in $colsToDelete) {
    $found = $ws.Range("1:1").Find($colName, $null, $null, 1)
    if ($found) {
        $found.EntireColumn.Delete()
    }
}

# scientificNameAuthorship's sample/template value moves from
# ${summary.Author} to ${summary.authority}.
$authCol = $ws.Range("1:1").Find("scientificNameAuthorship", $null, $null, 1)
if ($authCol) {
    $ws.Cells.Item(2, $authCol.Column).Value = "`${summary.authority}"
}
